$d = $word.ActiveDocument

# --- Step 1: Rewrite paragraph 3 (main intro paragraph), splitting it into the
#     updated intro paragraph + the new "At this report" paragraph + one blank
#     paragraph, all inheriting the ind/firstLine=708 formatting of the
#     original paragraph. ---
$p3 = $d.Paragraphs.Item(3)
$p3Range = $p3.Range
$p3Content = $d.Range($p3Range.Start, $p3Range.End - 1)

$para3Text = "In this project, AC to DC motor driving systems are investigated. The input is given as AC voltages by using VARIAC and it is supposed that the motor is controlled by externally without changing the input AC voltages.  There are some topologies such as 3-Phase Thyristor Rectifier, 1-Phase Thyristor Rectifier, Diode Rectifier with Buck Converter to operates for this aim.  We, The Third Harmonics, chose the topology that depends on diode rectifiers with PWM control"
$para4Text = "At this report, the topology will be examined deeply and argued in respect to advantages and disadvantages.  Computer simulation will be shown part by part by using Simulink. The component selection of the project will be made as result of the computer simulation. Then, the prototype will be created and tested. "

$p3Content.Text = $para3Text + "`r" + $para4Text + "`r"

# --- Step 2: Add a second blank paragraph (also ind firstLine=708) ---
$p5 = $d.Paragraphs.Item(5)
$p5Range = $p5.Range
$p5Content = $d.Range($p5Range.Start, $p5Range.End - 1)
$p5Content.Text = "`r"

# --- Step 3: Delete the "TOPOLOGY SELECTION" heading paragraph and the blank
#     bold paragraph that followed it ---
# Paragraph layout is now:
# 1 INTRODUCTION
# 2 (blank)
# 3 (intro text)
# 4 (At this report text)
# 5 (blank, ind firstLine=708)
# 6 (blank, ind firstLine=708)
# 7 (blank, ind firstLine=708)   <- pre-existing blank paragraph
# 8 TOPOLOGY SELECTION
# 9 (blank bold)
# 10 (tab + "At this project " + space + bookmark)
$pTopo = $d.Paragraphs.Item(8)
$pBlankBold = $d.Paragraphs.Item(9)
$delRange = $d.Range($pTopo.Range.Start, $pBlankBold.Range.End)
$delRange.Delete()

# --- Step 4: Remove the tab + "At this project " runs from the final
#     paragraph, keeping the trailing space run and the bookmark ---
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$searchRange = $pLast.Range.Duplicate()
$found = $searchRange.Find.Execute("At this project ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$toDelete = $d.Range($searchRange.Start - 1, $searchRange.End)
$toDelete.Delete()
